# Update the "End" column values (column D) for the BGM table on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 61875
$ws.Range("D4").Value = 42667
$ws.Range("D5").Value = 28872

# Move the active selection from D6 to D5
$ws.Range("D5").Select()
